# Apply data updates to the "Inscricoes" table per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E2").Value = 104
$ws.Range("E3").Value = 44
$ws.Range("F3").Value = 31
$ws.Range("H3").Value = 32
$ws.Range("E5").Value = 146
$ws.Range("F5").Value = 101
$ws.Range("H5").Value = 112
$ws.Range("E7").Value = 35
$ws.Range("F7").Value = 23
$ws.Range("H7").Value = 27
$ws.Range("F9").Value = 6
$ws.Range("H9").Value = 7
$ws.Range("E10").Value = 626
$ws.Range("F10").Value = 338
$ws.Range("H10").Value = 434
$ws.Range("E11").Value = 413
$ws.Range("F11").Value = 229
$ws.Range("H11").Value = 293
$ws.Range("E12").Value = 625
$ws.Range("F12").Value = 370
$ws.Range("H12").Value = 456
$ws.Range("F13").Value = 83
$ws.Range("H13").Value = 117
$ws.Range("F14").Value = 77
$ws.Range("H14").Value = 111
$ws.Range("F15").Value = 86
$ws.Range("H15").Value = 136
$ws.Range("F16").Value = 124
$ws.Range("H16").Value = 172
$ws.Range("E17").Value = 112
$ws.Range("F17").Value = 62
$ws.Range("H17").Value = 86
$ws.Range("F18").Value = 28
$ws.Range("H18").Value = 45
$ws.Range("E20").Value = 95
$ws.Range("F20").Value = 36
$ws.Range("H20").Value = 73
$ws.Range("E21").Value = 148
$ws.Range("F21").Value = 87
$ws.Range("H21").Value = 118
$ws.Range("F22").Value = 100
$ws.Range("H22").Value = 142
$ws.Range("F23").Value = 107
$ws.Range("H23").Value = 158
$ws.Range("E24").Value = 237
$ws.Range("F24").Value = 135
$ws.Range("H24").Value = 165
$ws.Range("E25").Value = 301
$ws.Range("F25").Value = 164
$ws.Range("H25").Value = 224
$ws.Range("E26").Value = 172
$ws.Range("F26").Value = 105
$ws.Range("H26").Value = 130
$ws.Range("E27").Value = 356
$ws.Range("F27").Value = 189
$ws.Range("H27").Value = 270
$ws.Range("E28").Value = 214
$ws.Range("F28").Value = 102
$ws.Range("H28").Value = 154
$ws.Range("E29").Value = 182
$ws.Range("F29").Value = 110
$ws.Range("H29").Value = 151
$ws.Range("F30").Value = 144
$ws.Range("H30").Value = 196
$ws.Range("E31").Value = 78
$ws.Range("F31").Value = 35
$ws.Range("H31").Value = 63
$ws.Range("E32").Value = 196
$ws.Range("F32").Value = 124
$ws.Range("H32").Value = 162
$ws.Range("E33").Value = 313
$ws.Range("F33").Value = 170
$ws.Range("H33").Value = 260
$ws.Range("E34").Value = 238
$ws.Range("F34").Value = 163
$ws.Range("H34").Value = 201
$ws.Range("E35").Value = 164
$ws.Range("F35").Value = 111
$ws.Range("H35").Value = 138
$ws.Range("E36").Value = 84
$ws.Range("F36").Value = 53
$ws.Range("H36").Value = 63
$ws.Range("E37").Value = 181
$ws.Range("F37").Value = 98
$ws.Range("H37").Value = 135
$ws.Range("E38").Value = 98
$ws.Range("F38").Value = 60
$ws.Range("H38").Value = 76
$ws.Range("F39").Value = 98
$ws.Range("H39").Value = 149
$ws.Range("E40").Value = 281
$ws.Range("F40").Value = 140
$ws.Range("H40").Value = 220
$ws.Range("E41").Value = 418
$ws.Range("F41").Value = 210
$ws.Range("H41").Value = 302
$ws.Range("E42").Value = 419
$ws.Range("F42").Value = 240
$ws.Range("H42").Value = 301
$ws.Range("F43").Value = 73
$ws.Range("H43").Value = 100
$ws.Range("E44").Value = 335
$ws.Range("F44").Value = 176
$ws.Range("H44").Value = 244
$ws.Range("F45").Value = 88
$ws.Range("H45").Value = 127
$ws.Range("E46").Value = 360
$ws.Range("F46").Value = 201
$ws.Range("H46").Value = 264
$ws.Range("E47").Value = 502
$ws.Range("F47").Value = 275
$ws.Range("H47").Value = 367
$ws.Range("E48").Value = 241
$ws.Range("F48").Value = 111
$ws.Range("H48").Value = 155
$ws.Range("E49").Value = 310
$ws.Range("F49").Value = 154
$ws.Range("H49").Value = 241
$ws.Range("E50").Value = 261
$ws.Range("F50").Value = 138
$ws.Range("H50").Value = 209
$ws.Range("F51").Value = 124
$ws.Range("H51").Value = 198
